$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the staffName2 (column B) values for rows 82-84 to mirror
# the staffName (column A) values, so the examiner modal can verify
# the examiner name against a second copy of the staff name.
$ws.Range("B82").Value2 = $ws.Range("A82").Value2
$ws.Range("B83").Value2 = $ws.Range("A83").Value2
$ws.Range("B84").Value2 = $ws.Range("A84").Value2

# Update the active window view/selection to match the saved state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("R64").Select()
